$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.01037833333333333
$ws.Range("H2").Value = 0.031135
$ws.Range("I2").Value = 0.02114284782989566
$ws.Range("J2").Value = 0.02114284782989566
$ws.Range("M2").Value = 5.575746
$ws.Range("N2").Value = 16.727238
$ws.Range("O2").Value = 0.069238947264747
$ws.Range("P2").Value = 0.069238947264747
$ws.Range("Q2").Value = 0.05786695057000001
$ws.Range("R2").Value = 0.5208025551300001
$ws.Range("S2").Value = 0.001463908525920716
$ws.Range("T2").Value = 0.001463908525920716
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.01037833333333333
$ws.Range("H3").Value = 0.031135
$ws.Range("I3").Value = 0.02114284782989566
$ws.Range("J3").Value = 0.02114284782989566
$ws.Range("O3").Value = 0.8150593598279631
$ws.Range("P3").Value = 0.815059359827963
$ws.Range("Q3").Value = 0.6811917504527779
$ws.Range("R3").Value = 6.130725754075001
$ws.Range("S3").Value = 0.0172326760171748
$ws.Range("T3").Value = 0.01723267601717479
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.01037833333333333
$ws.Range("H4").Value = 0.031135
$ws.Range("I4").Value = 0.02114284782989566
$ws.Range("J4").Value = 0.02114284782989566
$ws.Range("M4").Value = 9.317346333333333
$ws.Range("N4").Value = 27.952039
$ws.Range("O4").Value = 0.11570169290729
$ws.Range("P4").Value = 0.11570169290729
$ws.Range("Q4").Value = 0.09669852602944445
$ws.Range("R4").Value = 0.870286734265
$ws.Range("S4").Value = 0.00244626328680015
$ws.Range("T4").Value = 0.00244626328680015
$ws.Range("G5").Value = 0.4265683333333333
$ws.Range("I5").Value = 0.8690094132698448
$ws.Range("J5").Value = 0.8690094132698448
$ws.Range("M5").Value = 5.575746
$ws.Range("N5").Value = 16.727238
$ws.Range("O5").Value = 0.069238947264747
$ws.Range("P5").Value = 0.069238947264747
$ws.Range("Q5").Value = 2.378436678309999
$ws.Range("R5").Value = 21.40593010479
$ws.Range("S5").Value = 0.06016929693795951
$ws.Range("T5").Value = 0.06016929693795951
$ws.Range("G6").Value = 0.4265683333333333
$ws.Range("I6").Value = 0.8690094132698448
$ws.Range("J6").Value = 0.8690094132698448
$ws.Range("O6").Value = 0.8150593598279631
$ws.Range("P6").Value = 0.815059359827963
$ws.Range("Q6").Value = 27.99821708730277
$ws.Range("S6").Value = 0.7082942560641935
$ws.Range("T6").Value = 0.7082942560641934
$ws.Range("G7").Value = 0.4265683333333333
$ws.Range("I7").Value = 0.8690094132698448
$ws.Range("J7").Value = 0.8690094132698448
$ws.Range("M7").Value = 9.317346333333333
$ws.Range("N7").Value = 27.952039
$ws.Range("O7").Value = 0.11570169290729
$ws.Range("P7").Value = 0.11570169290729
$ws.Range("Q7").Value = 3.974484896499444
$ws.Range("R7").Value = 35.770364068495
$ws.Range("S7").Value = 0.1005458602676918
$ws.Range("T7").Value = 0.1005458602676918
$ws.Range("G8").Value = 0.05392066666666667
$ws.Range("H8").Value = 0.161762
$ws.Range("I8").Value = 0.1098477389002595
$ws.Range("J8").Value = 0.1098477389002595
$ws.Range("M8").Value = 5.575746
$ws.Range("N8").Value = 16.727238
$ws.Range("O8").Value = 0.069238947264747
$ws.Range("P8").Value = 0.069238947264747
$ws.Range("Q8").Value = 0.300647941484
$ws.Range("R8").Value = 2.705831473356
$ws.Range("S8").Value = 0.007605741800866768
$ws.Range("T8").Value = 0.007605741800866768
$ws.Range("G9").Value = 0.05392066666666667
$ws.Range("H9").Value = 0.161762
$ws.Range("I9").Value = 0.1098477389002595
$ws.Range("J9").Value = 0.1098477389002595
$ws.Range("O9").Value = 0.8150593598279631
$ws.Range("P9").Value = 0.815059359827963
$ws.Range("Q9").Value = 3.539134091432222
$ws.Range("R9").Value = 31.85220682289
$ws.Range("S9").Value = 0.08953242774659478
$ws.Range("T9").Value = 0.08953242774659477
$ws.Range("G10").Value = 0.05392066666666667
$ws.Range("H10").Value = 0.161762
$ws.Range("I10").Value = 0.1098477389002595
$ws.Range("J10").Value = 0.1098477389002595
$ws.Range("M10").Value = 9.317346333333333
$ws.Range("N10").Value = 27.952039
$ws.Range("O10").Value = 0.11570169290729
$ws.Range("P10").Value = 0.11570169290729
$ws.Range("Q10").Value = 0.5023975258575555
$ws.Range("R10").Value = 4.521577732718
$ws.Range("S10").Value = 0.012709569352798
$ws.Range("T10").Value = 0.012709569352798